$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Translate the prompt text in column B (rows 2-5) from English back to Chinese.
$ws.Range("B2").Value = "一间教室的典型布置"
$ws.Range("B3").Value = "窗户上雨滴顺着玻璃滑下"
$ws.Range("B4").Value = "天空中布满云朵"
$ws.Range("B5").Value = "火车停在火车站的情形"
